$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6743817574864011
$ws.Range("C2").Value = 0.2449791818670235
$ws.Range("D2").Value = 0.05732154878420204
$ws.Range("E2").Value = 0.1074546111888708
$ws.Range("F2").Value = 2.878803597020351
$ws.Range("I2").Value = 1.501337203754431
$ws.Range("J2").Value = 0.1807450889423734
$ws.Range("K2").Value = 0.9121141558723593
$ws.Range("M2").Value = 0.3575122264537711
$ws.Range("N2").Value = 2.967849150855329
$ws.Range("B3").Value = 0.640689483724401
$ws.Range("C3").Value = 0.2365572323610081
$ws.Range("D3").Value = 0.05632725425287077
$ws.Range("E3").Value = 0.1067583391954621
$ws.Range("F3").Value = 2.871212034085929
$ws.Range("I3").Value = 1.50149608696001
$ws.Range("J3").Value = 0.1802297368062185
$ws.Range("K3").Value = 0.8708599391783878
$ws.Range("M3").Value = 0.3475142182274169
$ws.Range("N3").Value = 2.984003695558819
$ws.Range("B4").Value = 0.6203843842383776
$ws.Range("C4").Value = 0.2315308718671645
$ws.Range("D4").Value = 0.05574239339185283
$ws.Range("E4").Value = 0.1063766714217707
$ws.Range("F4").Value = 2.867788914015279
$ws.Range("I4").Value = 1.502191119430293
$ws.Range("J4").Value = 0.1799882657655303
$ws.Range("K4").Value = 0.8460555906855518
$ws.Range("M4").Value = 0.341572948509139
$ws.Range("N4").Value = 2.994686754515527
$ws.Range("B5").Value = 0.6122060924509753
$ws.Range("C5").Value = 0.2295189391075496
$ws.Range("D5").Value = 0.05551053855027277
$ws.Range("E5").Value = 0.1062326869378829
$ws.Range("F5").Value = 2.866705290557562
$ws.Range("I5").Value = 1.502624559269861
$ws.Range("J5").Value = 0.1799087227176059
$ws.Range("K5").Value = 0.8360799160352883
$ws.Range("M5").Value = 0.3392015573179492
$ws.Range("N5").Value = 2.999232255595537
$ws.Range("B6").Value = 0.610853907305227
$ws.Range("C6").Value = 0.2291870532635727
$ws.Range("D6").Value = 0.05547243163456983
$ws.Range("E6").Value = 0.1062094764981758
$ws.Range("F6").Value = 2.866544158805524
$ws.Range("I6").Value = 1.502705602965705
$ws.Range("J6").Value = 0.1798966539491644
$ws.Range("K6").Value = 0.8344314551528953
$ws.Range("M6").Value = 0.3388107943482268
$ws.Range("N6").Value = 2.999998630959212
$ws.Range("B7").Value = 0.6202736992592861
$ws.Range("C7").Value = 0.2315035911230154
$ws.Range("D7").Value = 0.05573924023145338
$ws.Range("E7").Value = 0.1063746828115804
$ws.Range("F7").Value = 2.867773039368018
$ws.Range("I7").Value = 1.50219635681114
$ws.Range("J7").Value = 0.1799871166501745
$ws.Range("K7").Value = 0.845920519500396
$ws.Range("M7").Value = 0.3415407656991931
$ws.Range("N7").Value = 2.994747279165175
$ws.Range("B8").Value = 0.6626854175413257
$ws.Range("C8").Value = 0.2420452178231756
$ws.Range("D8").Value = 0.05697341392967559
$ws.Range("E8").Value = 0.1072050326718959
$ws.Range("F8").Value = 2.875929034661112
$ws.Range("I8").Value = 1.501267957186904
$ws.Range("J8").Value = 0.1805518433626148
$ws.Range("K8").Value = 0.8977805665754204
$ws.Range("M8").Value = 0.3540239374602336
$ws.Range("N8").Value = 2.973260604696357
$ws.Range("B9").Value = 0.7488863686125455
$ws.Range("C9").Value = 0.2638694768930918
$ws.Range("D9").Value = 0.05959579168368379
$ws.Range("E9").Value = 0.109196444780352
$ws.Range("F9").Value = 2.901751739211321
$ws.Range("I9").Value = 1.504191089102676
$ws.Range("J9").Value = 0.1822538565913447
$ws.Range("K9").Value = 1.00365502748204
$ws.Range("M9").Value = 0.3800704623262945
$ws.Range("N9").Value = 2.937191387577698
$ws.Range("B10").Value = 0.8140741764142945
$ws.Range("C10").Value = 0.2806132397461738
$ws.Range("D10").Value = 0.06164426381064203
$ws.Range("E10").Value = 0.1108802921017933
$ws.Range("F10").Value = 2.926728430311385
$ws.Range("I10").Value = 1.509237025052016
$ws.Range("J10").Value = 0.1838669668880968
$ws.Range("K10").Value = 1.084002982979115
$ws.Range("M10").Value = 0.4001642874106679
$ws.Range("N10").Value = 2.914393943654261
$ws.Range("B11").Value = 0.8441349944173737
$ws.Range("C11").Value = 0.2883861370700345
$ws.Range("D11").Value = 0.06260233349875932
$ws.Range("E11").Value = 0.1116941493018047
$ws.Range("F11").Value = 2.939398102018586
$ws.Range("I11").Value = 1.512163465637791
$ws.Range("J11").Value = 0.1846796235980008
$ws.Range("K11").Value = 1.121115652017352
$ws.Range("M11").Value = 0.4095140276668587
$ws.Range("N11").Value = 2.904827861580458
$ws.Range("B12").Value = 0.8555767207441818
$ws.Range("C12").Value = 0.2913520679206272
$ws.Range("D12").Value = 0.06296886943026436
$ws.Range("E12").Value = 0.1120092046633445
$ws.Range("F12").Value = 2.944383975270327
$ws.Range("I12").Value = 1.51336246598725
$ws.Range("J12").Value = 0.1849986920390023
$ws.Range("K12").Value = 1.135250188445468
$ws.Range("M12").Value = 0.4130845771539242
$ws.Range("N12").Value = 2.901321246922592
$ws.Range("B13").Value = 0.8531099470739605
$ws.Range("C13").Value = 0.2907122998383613
$ws.Range("D13").Value = 0.06288976373709687
$ws.Range("E13").Value = 0.1119410468003359
$ws.Range("F13").Value = 2.943301808790125
$ws.Range("I13").Value = 1.513100199210569
$ws.Range("J13").Value = 0.1849294709768188
$ws.Range("K13").Value = 1.132202471241527
$ws.Range("M13").Value = 0.4123142615817486
$ws.Range("N13").Value = 2.902071304036369
$ws.Range("B14").Value = 0.8450751434208144
$ws.Range("C14").Value = 0.2886296946395817
$ws.Range("D14").Value = 0.06263241394911745
$ws.Range("E14").Value = 0.1117199316067214
$ws.Range("F14").Value = 2.939804521514901
$ws.Range("I14").Value = 1.512260287461295
$ws.Range("J14").Value = 0.1847056464470143
$ws.Range("K14").Value = 1.122276890207274
$ws.Range("M14").Value = 0.4098071777994647
$ws.Range("N14").Value = 2.904537047519071
$ws.Range("B15").Value = 0.8401611892977883
$ws.Range("C15").Value = 0.2873569712158712
$ws.Range("D15").Value = 0.06247526529736547
$ws.Range("E15").Value = 0.1115853857047071
$ws.Range("F15").Value = 2.937686837126563
$ws.Range("I15").Value = 1.511757647259635
$ws.Range("J15").Value = 0.1845700232114496
$ws.Range("K15").Value = 1.11620770333667
$ws.Range("M15").Value = 0.4082754229960628
$ws.Range("N15").Value = 2.90606247872131
$ws.Range("B16").Value = 0.8121178022982463
$ws.Range("C16").Value = 0.280108408354522
$ws.Range("D16").Value = 0.06158217626489915
$ws.Range("E16").Value = 0.1108280662873966
$ws.Range("F16").Value = 2.925926759590197
$ws.Range("I16").Value = 1.509058479864869
$ws.Range("J16").Value = 0.1838154440020503
$ws.Range("K16").Value = 1.081588887932213
$ws.Range("M16").Value = 0.399557461040942
$ws.Range("N16").Value = 2.915035313083251
$ws.Range("B17").Value = 0.7950181387379587
$ws.Range("C17").Value = 0.2757016667407868
$ws.Range("D17").Value = 0.06104098427140059
$ws.Range("E17").Value = 0.110375722372531
$ws.Range("F17").Value = 2.919047326553269
$ws.Range("I17").Value = 1.507564304104683
$ws.Range("J17").Value = 0.1833727246172998
$ws.Range("K17").Value = 1.06049527454519
$ws.Range("M17").Value = 0.3942627626384478
$ws.Range("N17").Value = 2.92074603930115
$ws.Range("B18").Value = 0.7852211383545011
$ws.Range("C18").Value = 0.2731817124447389
$ws.Range("D18").Value = 0.06073217387625363
$ws.Range("E18").Value = 0.1101200526046746
$ws.Range("F18").Value = 2.91521353415591
$ws.Range("I18").Value = 1.50676428320196
$ws.Range("J18").Value = 0.1831255056460535
$ws.Range("K18").Value = 1.048415684394428
$ws.Range("M18").Value = 0.3912370627674093
$ws.Range("N18").Value = 2.924106420233088
$ws.Range("B19").Value = 0.7819106194916969
$ws.Range("C19").Value = 0.27233101990862
$ws.Range("D19").Value = 0.06062804103575559
$ws.Range("E19").Value = 0.1100342617550112
$ws.Range("F19").Value = 2.913936613159876
$ws.Range("I19").Value = 1.506503607648845
$ws.Range("J19").Value = 0.1830430763521562
$ws.Range("K19").Value = 1.044334829292382
$ws.Range("M19").Value = 0.3902159925305355
$ws.Range("N19").Value = 2.925257188238902
$ws.Range("B20").Value = 0.7968344667630731
$ws.Range("C20").Value = 0.2761692514375227
$ws.Range("D20").Value = 0.06109833980038104
$ws.Range("E20").Value = 0.1104234088802087
$ws.Range("F20").Value = 2.919766915333184
$ws.Range("I20").Value = 1.507717214491407
$ws.Range("J20").Value = 0.1834190847892572
$ws.Range("K20").Value = 1.062735251206249
$ws.Range("M20").Value = 0.3948243567975851
$ws.Range("N20").Value = 2.92013028466063
$ws.Range("B21").Value = 0.847433577407287
$ws.Range("C21").Value = 0.2892407950537006
$ws.Range("D21").Value = 0.06270790269692839
$ws.Range("E21").Value = 0.1117846923179648
$ws.Range("F21").Value = 2.940826652537268
$ws.Range("I21").Value = 1.512504524392519
$ws.Range("J21").Value = 0.1847710815999406
$ws.Range("K21").Value = 1.125190082055241
$ws.Range("M21").Value = 0.4105427550615914
$ws.Range("N21").Value = 2.903809653421192
$ws.Range("B22").Value = 0.8808430603188242
$ws.Range("C22").Value = 0.2979149996755268
$ws.Range("D22").Value = 0.06378160459419746
$ws.Range("E22").Value = 0.1127143778317254
$ws.Range("F22").Value = 2.955687105282493
$ws.Range("I22").Value = 1.516162720722164
$ws.Range("J22").Value = 0.1857207449620475
$ws.Range("K22").Value = 1.166478722782898
$ws.Range("M22").Value = 0.4209905147919315
$ws.Range("N22").Value = 2.893818531908096
$ws.Range("B23").Value = 0.8629807143674384
$ws.Range("C23").Value = 0.2932733889663837
$ws.Range("D23").Value = 0.06320656957997528
$ws.Range("E23").Value = 0.112214532281012
$ws.Range("F23").Value = 2.947655418055348
$ws.Range("I23").Value = 1.514161803800796
$ws.Range("J23").Value = 0.1852078491926861
$ws.Range("K23").Value = 1.144399131690022
$ws.Range("M23").Value = 0.4153983611092684
$ws.Range("N23").Value = 2.899089134562857
$ws.Range("B24").Value = 0.7960131994963149
$ws.Range("C24").Value = 0.2759578142245687
$ws.Range("D24").Value = 0.06107240211227349
$ws.Range("E24").Value = 0.1104018361386139
$ws.Range("F24").Value = 2.919441211440756
$ws.Range("I24").Value = 1.507647899917096
$ws.Range("J24").Value = 0.1833981025972804
$ws.Range("K24").Value = 1.061722410081188
$ws.Range("M24").Value = 0.3945704030956492
$ws.Range("N24").Value = 2.920408426820728
$ws.Range("B25").Value = 0.725241371014846
$ws.Range("C25").Value = 0.2578413610751227
$ws.Range("D25").Value = 0.05886487493746984
$ws.Range("E25").Value = 0.1086189053800588
$ws.Range("F25").Value = 2.893712463016044
$ws.Range("I25").Value = 1.502891803042687
$ws.Range("J25").Value = 0.1817297464410146
$ws.Range("K25").Value = 0.974564404995192
$ws.Range("M25").Value = 0.3728562075693205
$ws.Range("N25").Value = 2.946299194559316
